$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "plan"
$ws.Cells.Item(1,2).Value = "geo_return_23yr"
$ws.Cells.Item(1,3).Value = "mva_2023"

# Data rows (plan name, geo_return_23yr, mva_2023)
$data = @(
  @("Utah Retirement Systems, Contributory Retirement System", 0.115117531699273, 1226842000),
  @("Utah Retirement Systems, Firefighters Retirement System", 0.115117531699273, 1815120000),
  @("Utah Retirement Systems, Tier 2 Public Employees Contributory Retirement System", 0.102213448486189, 1673551000),
  @("Minnesota State Retirement System (MSRS) Judges Retirement Plan", 0.100837406897798, 268987000),
  @("Minnesota State Retirement System (MSRS) State Patrol Retirement Plan", 0.0928089303297643, 943099000),
  @("Minnesota Public Employees Retirement Association (MPERA) Local Government Correctional Employees Plan", 0.0919642003464134, 1067200000),
  @("Nebraska County Employees Retirement System", 0.0882012942465855, 664891026),
  @("Kansas Public Employees' Retirement System", 0.0857686523972012, 25800659628),
  @("Tri-County Metro Of Oregon Bargaining Unit Defined Benefit Plan", 0.0851986968968603, 708822000),
  @("Arkansas Judicial Retirement Plan", 0.0832043607418855, 308094805)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
